# "updated url to app"
#
# The bullet paragraph that describes the streamlit deployment used to be a
# single run of text. We split it into three runs so that the words
# "the app " become a hyperlink pointing at the deployed streamlit app,
# while the rest of the sentence keeps its original (non-linked) formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The big bullet textbox ("TextBox 3") holds the paragraph we need to edit.
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

# Locate the paragraph that mentions the streamlit deployment instead of
# hard-coding its index, so the script keeps working even if other bullets
# shift around.
$paraCount = $tr.Paragraphs().Count
$targetParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $paraText = $tr.Paragraphs($i).Text
    if ($paraText.IndexOf("For your convinience") -ge 0) {
        $targetParaIndex = $i
        break
    }
}

$para = $tr.Paragraphs($targetParaIndex)
$fullText = $para.Text

# Figure out where "the app " sits inside the paragraph so we can carve it
# into its own run without disturbing the surrounding text/formatting.
$linkText = "the app "
$linkStart0 = $fullText.IndexOf($linkText)

# Characters() is 1-based and relative to the paragraph itself.
$linkStart = $linkStart0 + 1
$linkLength = $linkText.Length

$linkRun = $para.Characters($linkStart, $linkLength)

# Turning this run into a hyperlink automatically splits the paragraph's
# single run into three runs (before / link / after), matching the target
# OOXML (a new <a:hlinkClick r:id="rIdX"/> run plus a fresh relationship).
$linkRun.ActionSettings(1).Hyperlink.Address = "https://llm-for-qa.streamlit.app/"

Write-Host "Linked text:" $linkRun.Text
Write-Host "Paragraph now:" $para.Text
